$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, [string]$value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.ClearFormats()
}

Set-TextValue $ws.Range("D2") "64.199.37"
Set-TextValue $ws.Range("E2") "  -2.33%  "
Set-TextValue $ws.Range("D3") "3.176.53"
Set-TextValue $ws.Range("E3") "  -7.69%  "
Set-TextValue $ws.Range("E4") "  +0.07%  "
Set-TextValue $ws.Range("D5") "556.07"
Set-TextValue $ws.Range("E5") "  -4.58%  "
Set-TextValue $ws.Range("D6") "170.42"
Set-TextValue $ws.Range("E6") "  -1.25%  "
Set-TextValue $ws.Range("E7") "  +0.05%  "
Set-TextValue $ws.Range("D8") "0.600"
Set-TextValue $ws.Range("E8") "  +0.57%  "
Set-TextValue $ws.Range("D9") "3.175.69"
Set-TextValue $ws.Range("E9") "  -7.63%  "
Set-TextValue $ws.Range("E10") "  -5.73%  "
Set-TextValue $ws.Range("D11") "6.60"
Set-TextValue $ws.Range("E11") "  -4.09%  "
Set-TextValue $ws.Range("D12") "0.394"
Set-TextValue $ws.Range("E12") "  -3.53%  "
Set-TextValue $ws.Range("D13") "3.731.57"
Set-TextValue $ws.Range("E13") "  -7.54%  "
Set-TextValue $ws.Range("E14") "  -0.27%  "
Set-TextValue $ws.Range("D15") "27.27"
Set-TextValue $ws.Range("E15") "  -5.10%  "
Set-TextValue $ws.Range("D16") "64.402.47"
Set-TextValue $ws.Range("E16") "  -2.14%  "
Set-TextValue $ws.Range("D17") "0.0000161"
Set-TextValue $ws.Range("E17") "  -5.33%  "
Set-TextValue $ws.Range("D18") "3.193.34"
Set-TextValue $ws.Range("E18") "  -7.10%  "
Set-TextValue $ws.Range("D19") "5.61"
Set-TextValue $ws.Range("E19") "  -5.16%  "
Set-TextValue $ws.Range("D20") "12.93"
Set-TextValue $ws.Range("E20") "  -6.54%  "
Set-TextValue $ws.Range("D21") "351.50"
Set-TextValue $ws.Range("E21") "  -3.84%  "
Set-TextValue $ws.Range("D22") "7.15"
Set-TextValue $ws.Range("E22") "  -6.63%  "
Set-TextValue $ws.Range("E23") "  +0.12%  "
Set-TextValue $ws.Range("D24") "68.62"
Set-TextValue $ws.Range("E24") "  -5.65%  "
Set-TextValue $ws.Range("E25") "  -3.22%  "
Set-TextValue $ws.Range("D26") "0.499"
Set-TextValue $ws.Range("E26") "  -6.34%  "
Set-TextValue $ws.Range("D27") "9.37"
Set-TextValue $ws.Range("E27") "  -3.98%  "
Set-TextValue $ws.Range("D28") "0.174"
Set-TextValue $ws.Range("E28") "  -1.94%  "
Set-TextValue $ws.Range("E29") "  -0.18%  "
Set-TextValue $ws.Range("B30") "USDe"
Set-TextValue $ws.Range("C30") "https://coinranking.com/coin/exbfr2U-0+usde-usde"
Set-TextValue $ws.Range("D30") "1.00"
Set-TextValue $ws.Range("E30") "  -0.01%  "
Set-TextValue $ws.Range("B31") "NEARProtocol"
Set-TextValue $ws.Range("C31") "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
Set-TextValue $ws.Range("D31") "5.58"
Set-TextValue $ws.Range("E31") "  -1.99%  "
Set-TextValue $ws.Range("D32") "1.88"
Set-TextValue $ws.Range("E32") "  -4.89%  "
Set-TextValue $ws.Range("D33") "22.05"
Set-TextValue $ws.Range("E33") "  -7.02%  "
Set-TextValue $ws.Range("D34") "6.57"
Set-TextValue $ws.Range("E34") "  -6.44%  "
Set-TextValue $ws.Range("D35") "1.18"
Set-TextValue $ws.Range("E35") "  -8.75%  "
Set-TextValue $ws.Range("D36") "158.04"
Set-TextValue $ws.Range("E36") "  -2.06%  "
Set-TextValue $ws.Range("E37") "  -6.22%  "
Set-TextValue $ws.Range("D38") "0.806"
Set-TextValue $ws.Range("E38") "  -8.54%  "
Set-TextValue $ws.Range("D39") "26.23"
Set-TextValue $ws.Range("E39") "  -9.61%  "
Set-TextValue $ws.Range("D40") "2.50"
Set-TextValue $ws.Range("E40") "  -3.75%  "
Set-TextValue $ws.Range("B41") "Stacks"
Set-TextValue $ws.Range("C41") "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
Set-TextValue $ws.Range("D41") "1.67"
Set-TextValue $ws.Range("E41") "  -4.13%  "
Set-TextValue $ws.Range("B42") "Maker"
Set-TextValue $ws.Range("C42") "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
Set-TextValue $ws.Range("D42") "2.640.00"
Set-TextValue $ws.Range("E42") "  -4.41%  "
Set-TextValue $ws.Range("D43") "6.02"
Set-TextValue $ws.Range("E43") "  -6.22%  "
Set-TextValue $ws.Range("D44") "4.12"
Set-TextValue $ws.Range("E44") "  -7.07%  "
Set-TextValue $ws.Range("D45") "0.0649"
Set-TextValue $ws.Range("E45") "  -4.48%  "
Set-TextValue $ws.Range("D46") "38.76"
Set-TextValue $ws.Range("E46") "  -3.08%  "
Set-TextValue $ws.Range("D47") "319.49"
Set-TextValue $ws.Range("E47") "  -0.75%  "
Set-TextValue $ws.Range("D48") "23.18"
Set-TextValue $ws.Range("E48") "  -3.83%  "
Set-TextValue $ws.Range("D49") "0.0268"
Set-TextValue $ws.Range("E49") "  -7.02%  "
Set-TextValue $ws.Range("D50") "0.101"
Set-TextValue $ws.Range("E50") "  +0.06%  "
Set-TextValue $ws.Range("E51") "  +0.06%  "
